$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 558.1818
$ws.Range("I28").Value = 381.55554
$ws.Range("J28").Value = 1353
$ws.Range("K28").Value = 381.55554
$ws.Range("L28").Value = 1353
$ws.Range("M28").Value = 103.44446
$ws.Range("N28").Value = -2323

$ws.Range("H62").Value = 50002760
$ws.Range("J62").Value = 2832.6667
$ws.Range("L62").Value = 2832.6667
$ws.Range("N62").Value = -4080.6667

$ws.Range("H65").Value = 50002760
$ws.Range("J65").Value = 2832.6667
$ws.Range("L65").Value = 14163.3335
$ws.Range("N65").Value = -20403.3335

$ws.Range("H103").Value = 1000.9231
$ws.Range("I103").Value = 300
$ws.Range("J103").Value = 1312.4445
$ws.Range("K103").Value = 900
$ws.Range("L103").Value = 3937.3335
$ws.Range("M103").Value = -314
$ws.Range("N103").Value = -5109.333500000001

$ws.Range("H116").Value = 12823345
$ws.Range("I116").Value = 15387414
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 15387414
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -15383972
$ws.Range("N116").Value = -9884

$ws.Range("H118").Value = 520
$ws.Range("I118").Value = 587.5
$ws.Range("J118").Value = 250
$ws.Range("K118").Value = 1762.5
$ws.Range("L118").Value = 750
$ws.Range("M118").Value = -105.5
$ws.Range("N118").Value = -4064

$ws.Range("H129").Value = 981.77014
$ws.Range("I129").Value = 758.6
$ws.Range("J129").Value = 995.37805
$ws.Range("K129").Value = 2275.8
$ws.Range("L129").Value = 2986.13415
$ws.Range("M129").Value = 2724.2
$ws.Range("N129").Value = -12986.13415

$ws.Range("H131").Value = 65692.69
$ws.Range("I131").Value = 144652.86
$ws.Range("J131").Value = 4279.222
$ws.Range("K131").Value = 433958.58
$ws.Range("L131").Value = 12837.666
$ws.Range("M131").Value = -428918.58
$ws.Range("N131").Value = -22917.666

$ws.Range("H137").Value = 7172883
$ws.Range("I137").Value = 10417550
$ws.Range("J137").Value = 93608.82000000001
$ws.Range("K137").Value = 31252650
$ws.Range("L137").Value = 280826.46
$ws.Range("M137").Value = -31250100
$ws.Range("N137").Value = -285926.46

$ws.Range("H139").Value = 11593.333
$ws.Range("J139").Value = 11593.333
$ws.Range("L139").Value = 11593.333
$ws.Range("N139").Value = -21873.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1991.8889
$ws.Range("I2").Value = 2139.3076
$ws.Range("J2").Value = 1608.6
$ws.Range("K2").Value = 2139.3076
$ws.Range("L2").Value = 1608.6
$ws.Range("M2").Value = -2026.3076
$ws.Range("N2").Value = -1834.6

$ws.Range("H110").Value = 939.2632
$ws.Range("I110").Value = 528.8333
$ws.Range("J110").Value = 1642.8572
$ws.Range("K110").Value = 528.8333
$ws.Range("L110").Value = 1642.8572
$ws.Range("M110").Value = 1516.1667
$ws.Range("N110").Value = -5732.8572

$ws.Range("H116").Value = 1991.8889
$ws.Range("I116").Value = 2139.3076
$ws.Range("J116").Value = 1608.6
$ws.Range("K116").Value = 2139.3076
$ws.Range("L116").Value = 1608.6
$ws.Range("M116").Value = 154.6923999999999
$ws.Range("N116").Value = -6196.6

$ws.Range("H122").Value = 1421.0667
$ws.Range("I122").Value = 1280.5555
$ws.Range("J122").Value = 1631.8334
$ws.Range("K122").Value = 3841.6665
$ws.Range("L122").Value = 4895.5002
$ws.Range("M122").Value = -1391.6665
$ws.Range("N122").Value = -9795.5002

$ws.Range("H132").Value = 3158.8096
$ws.Range("I132").Value = 2760.4285
$ws.Range("J132").Value = 3955.5715
$ws.Range("K132").Value = 8281.2855
$ws.Range("L132").Value = 11866.7145
$ws.Range("M132").Value = -5751.2855
$ws.Range("N132").Value = -16926.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1991.8889
$ws.Range("I3").Value = 2139.3076
$ws.Range("J3").Value = 1608.6
$ws.Range("K3").Value = 2139.3076
$ws.Range("L3").Value = 1608.6
$ws.Range("M3").Value = -2025.3076
$ws.Range("N3").Value = -1836.6

$ws.Range("H99").Value = 2834
$ws.Range("I99").Value = 1604
$ws.Range("J99").Value = 3602.75
$ws.Range("K99").Value = 1604
$ws.Range("L99").Value = 3602.75
$ws.Range("M99").Value = -106
$ws.Range("N99").Value = -6598.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3455.5557
$ws.Range("I16").Value = 2250
$ws.Range("J16").Value = 5866.6665
$ws.Range("K16").Value = 2250
$ws.Range("L16").Value = 5866.6665
$ws.Range("M16").Value = -1963
$ws.Range("N16").Value = -6440.6665

$ws.Range("H31").Value = 2674.848
$ws.Range("I31").Value = 2003.8857
$ws.Range("K31").Value = 2003.8857
$ws.Range("M31").Value = -1708.8857

$ws.Range("H34").Value = 2674.848
$ws.Range("I34").Value = 2003.8857
$ws.Range("K34").Value = 2003.8857
$ws.Range("M34").Value = -1801.8857

$ws.Range("H113").Value = 3455.5557
$ws.Range("I113").Value = 2250
$ws.Range("J113").Value = 5866.6665
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 5866.6665
$ws.Range("M113").Value = -80
$ws.Range("N113").Value = -10206.6665

$ws.Range("H122").Value = 1373.2333
$ws.Range("I122").Value = 1355.4117
$ws.Range("J122").Value = 1396.5385
$ws.Range("K122").Value = 4066.2351
$ws.Range("L122").Value = 4189.6155
$ws.Range("M122").Value = -1616.2351
$ws.Range("N122").Value = -9089.6155

$ws.Range("H134").Value = 6044.375
$ws.Range("I134").Value = 1392.6666
$ws.Range("J134").Value = 19999.5
$ws.Range("K134").Value = 4177.9998
$ws.Range("L134").Value = 59998.5
$ws.Range("M134").Value = -1642.9998
$ws.Range("N134").Value = -65068.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2155658.8
$ws.Range("I113").Value = 6896888
$ws.Range("J113").Value = 554.4545000000001
$ws.Range("K113").Value = 20690664
$ws.Range("L113").Value = 1663.3635
$ws.Range("M113").Value = -20688494
$ws.Range("N113").Value = -6003.3635

$ws.Range("H131").Value = 738.0417
$ws.Range("J131").Value = 904.5625
$ws.Range("L131").Value = 2713.6875
$ws.Range("N131").Value = -12793.6875

$ws.Range("H132").Value = 723001.4399999999
$ws.Range("I132").Value = 1214
$ws.Range("J132").Value = 1123994.5
$ws.Range("K132").Value = 10926
$ws.Range("L132").Value = 10115950.5
$ws.Range("M132").Value = -8396
$ws.Range("N132").Value = -10121010.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 590
$ws.Range("I107").Value = 590
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 590
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1330
$ws.Range("N107").Value = ""

$ws.Range("H113").Value = 6962.952
$ws.Range("I113").Value = 1859.5
$ws.Range("J113").Value = 11602.454
$ws.Range("K113").Value = 1859.5
$ws.Range("L113").Value = 11602.454
$ws.Range("M113").Value = 310.5
$ws.Range("N113").Value = -15942.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 8274.529
$ws.Range("I61").Value = 10447.167
$ws.Range("J61").Value = 3060.2
$ws.Range("K61").Value = 10447.167
$ws.Range("L61").Value = 3060.2
$ws.Range("M61").Value = -10245.167
$ws.Range("N61").Value = -3464.2

$ws.Range("H113").Value = 8274.529
$ws.Range("I113").Value = 10447.167
$ws.Range("J113").Value = 3060.2
$ws.Range("K113").Value = 10447.167
$ws.Range("L113").Value = 3060.2
$ws.Range("M113").Value = -8277.166999999999
$ws.Range("N113").Value = -7400.2

$ws.Range("H122").Value = 2869
$ws.Range("I122").Value = 2803.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8410.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5960.5
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1513.0667
$ws.Range("I132").Value = 782.8857400000001
$ws.Range("J132").Value = 4068.7
$ws.Range("K132").Value = 2348.65722
$ws.Range("L132").Value = 12206.1
$ws.Range("M132").Value = 181.3427799999999
$ws.Range("N132").Value = -17266.1
